$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the FilesTab Neo4j query text in B4: drop the File Type and Breed coalesce lines
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$ws.Range("B4").Value2 = $newFilesQuery

# Row 4 shrinks now that two lines were removed from the wrapped text
$ws.Rows.Item(4).RowHeight = 217.5

# Move the active selection to B4 (was C4)
$ws.Activate()
$ws.Range("B4").Select()
